$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 8, shifting existing rows 8-31 down to 9-32
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with the inserted data
$ws.Cells.Item(8, 1).Value = 3
$ws.Cells.Item(8, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(8, 3).Value = "Coquimbo"
$ws.Cells.Item(8, 4).Value = 44487
$ws.Cells.Item(8, 5).Value = 5
$ws.Cells.Item(8, 6).Value = 100112022
$ws.Cells.Item(8, 7).Value = "Arveja Verde"
$ws.Cells.Item(8, 8).Value = "Perfection"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 73
$ws.Cells.Item(8, 11).Value = 20000
$ws.Cells.Item(8, 12).Value = 21000
$ws.Cells.Item(8, 13).Value = 20521
$ws.Cells.Item(8, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(8, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(8, 16).Value = 821
$ws.Cells.Item(8, 17).Value = 25
$ws.Cells.Item(8, 18).Value = "Hortaliza"
